# Update scripts with new TPM values.
#
# The source data dropped the "ECs" sending-cluster row entirely, leaving
# only the "MuSCs" row (which moves up to row 2). Because the specificity
# metrics are computed relative to all rows present, with only a single
# remaining row the derived-specificity columns (I, J, S, T) normalize to 1,
# and a couple of values were re-emitted with cleaner floating point
# precision (H, M).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "ECs" row (original row 2); this shifts the "MuSCs" row up to
# row 2 and drops the now-unused "ECs" shared string automatically.
$ws.Rows.Item(2).Delete()

# Refresh the recalculated values on the remaining (former "MuSCs") row.
$ws.Range("G2").Value = 0.011155
$ws.Range("H2").Value = 0.033465
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("M2").Value = 0.038688
$ws.Range("Q2").Value = 0.00043156464
$ws.Range("R2").Value = 0.00388408176
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1
